$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Data update -----------------------------------------------------
# Fix typo in the "Unsupervised Machine  Learning" (double space) entry
# so it reads "Unsupervised Machine Learning" (single space). The "UML"
# short label in column B stays the same.
$ws.Range("A13").Value = "Unsupervised Machine Learning"
$ws.Range("B13").Value = "UML"

# Add two new rows of data at the bottom of the table (Haskell, and
# Visual Studio / VS).
$ws.Range("A44").Value = "Haskell"
$ws.Range("B44").Value = "Haskell"

$ws.Range("A45").Value = "Visual Studio"
$ws.Range("B45").Value = "VS"

# --- Cell style name fix ----------------------------------------------
# The workbook's built-in "Normal" cell style had been saved under its
# Polish localized name "Normalny"; rename it back to "Normal".
$wb.Styles.Item(1).Name = "Normal"

# --- View/selection update ---------------------------------------------
# Scroll down and select near the bottom of the newly extended table,
# matching where the sheet was left scrolled to after the edit.
$ws.Range("A43").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 43
$win.ScrollColumn = 1
$ws.Range("E48").Select()
